$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1701244813278008
$ws.Range("C2").Value = 0.5643153526970954
$ws.Range("J2").Value = 0.01659751037344398
$ws.Range("P2").Value = 0.1244813278008299
$ws.Range("S2").Value = 0.1244813278008299
$ws.Range("B3").Value = 0.007352941176470588
$ws.Range("J3").Value = 0.04411764705882353
$ws.Range("P3").Value = 0.7426470588235294
$ws.Range("S3").Value = 0.2058823529411765
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.7333333333333333
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.05194805194805195
$ws.Range("D6").Value = 0.01298701298701299
$ws.Range("E6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.06060606060606061
$ws.Range("J6").Value = 0.2424242424242424
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.1341991341991342
$ws.Range("R6").Value = 0.06060606060606061
$ws.Range("S6").Value = 0.4242424242424243
$ws.Range("B7").Value = 0.07053941908713693
$ws.Range("D7").Value = 0.01659751037344398
$ws.Range("E7").Value = 0.004149377593360996
$ws.Range("F7").Value = 0.04564315352697095
$ws.Range("J7").Value = 0.1452282157676349
$ws.Range("O7").Value = 0.008298755186721992
$ws.Range("Q7").Value = 0.1493775933609958
$ws.Range("R7").Value = 0.06639004149377593
$ws.Range("S7").Value = 0.4937759336099585
$ws.Range("B8").Value = 0.1075268817204301
$ws.Range("D8").Value = 0.02150537634408602
$ws.Range("F8").Value = 0.07168458781362007
$ws.Range("J8").Value = 0.1254480286738351
$ws.Range("O8").Value = 0.01792114695340502
$ws.Range("Q8").Value = 0.1577060931899641
$ws.Range("R8").Value = 0.07706093189964158
$ws.Range("S8").Value = 0.4211469534050179
$ws.Range("B9").Value = 0.1197916666666667
$ws.Range("D9").Value = 0.015625
$ws.Range("F9").Value = 0.08854166666666667
$ws.Range("J9").Value = 0.1145833333333333
$ws.Range("Q9").Value = 0.109375
$ws.Range("R9").Value = 0.08854166666666667
$ws.Range("S9").Value = 0.4635416666666667
$ws.Range("B10").Value = 0.09112333071484682
$ws.Range("D10").Value = 0.0180675569520817
$ws.Range("F10").Value = 0.05577376276512176
$ws.Range("J10").Value = 0.1241162608012569
$ws.Range("O10").Value = 0.00864100549882168
$ws.Range("Q10").Value = 0.2325216025137471
$ws.Range("R10").Value = 0.09112333071484682
$ws.Range("S10").Value = 0.3786331500392773
$ws.Range("G11").Value = 0.18125
$ws.Range("J11").Value = 0.08437500000000001
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5125
$ws.Range("S11").Value = 0.021875
$ws.Range("G12").Value = 0.7607655502392344
$ws.Range("J12").Value = 0.1626794258373206
$ws.Range("K12").Value = 0.009569377990430622
$ws.Range("L12").Value = 0.02870813397129187
$ws.Range("S12").Value = 0.03827751196172249
$ws.Range("G13").Value = 0.6274509803921569
$ws.Range("J13").Value = 0.2549019607843137
$ws.Range("S13").Value = 0.1176470588235294
$ws.Range("F15").Value = 0.02515723270440252
$ws.Range("H15").Value = 0.1823899371069182
$ws.Range("I15").Value = 0.08176100628930817
$ws.Range("J15").Value = 0.3396226415094339
$ws.Range("K15").Value = 0.06289308176100629
$ws.Range("M15").Value = 0.01257861635220126
$ws.Range("O15").Value = 0.0440251572327044
$ws.Range("S15").Value = 0.2515723270440252
$ws.Range("F16").Value = 0.0245398773006135
$ws.Range("H16").Value = 0.1717791411042945
$ws.Range("I16").Value = 0.0736196319018405
$ws.Range("J16").Value = 0.4110429447852761
$ws.Range("K16").Value = 0.1349693251533742
$ws.Range("M16").Value = 0.006134969325153374
$ws.Range("N16").Value = 0.006134969325153374
$ws.Range("O16").Value = 0.0245398773006135
$ws.Range("S16").Value = 0.147239263803681
$ws.Range("F17").Value = 0.01956521739130435
$ws.Range("H17").Value = 0.2260869565217391
$ws.Range("I17").Value = 0.09347826086956522
$ws.Range("J17").Value = 0.441304347826087
$ws.Range("K17").Value = 0.09347826086956522
$ws.Range("M17").Value = 0.01956521739130435
$ws.Range("O17").Value = 0.04347826086956522
$ws.Range("S17").Value = 0.06304347826086956
$ws.Range("F18").Value = 0.01951219512195122
$ws.Range("H18").Value = 0.2097560975609756
$ws.Range("I18").Value = 0.07317073170731707
$ws.Range("J18").Value = 0.424390243902439
$ws.Range("K18").Value = 0.1365853658536585
$ws.Range("M18").Value = 0.01951219512195122
$ws.Range("N18").Value = 0.00975609756097561
$ws.Range("O18").Value = 0.03902439024390244
$ws.Range("S18").Value = 0.06829268292682927
$ws.Range("F19").Value = 0.02809798270893372
$ws.Range("H19").Value = 0.2334293948126801
$ws.Range("I19").Value = 0.07420749279538905
$ws.Range("J19").Value = 0.345821325648415
$ws.Range("K19").Value = 0.09438040345821326
$ws.Range("M19").Value = 0.02449567723342939
$ws.Range("O19").Value = 0.04899135446685879
$ws.Range("S19").Value = 0.1505763688760807
